# Apply the 31-12-2022 13:49 UTC "Updated symbol list" crypto price/row refresh.
# Re-applies every changed cell from the commit diff (Coin/Link/Price/Volume columns,
# rows 2-50), including the row-shuffle in rows 10-14 and 41-42.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as *text* in the workbook (t="inlineStr").
# Assigning a bare numeric string (e.g. "246.72") to .Value lets Excel coerce it to a
# real number, which would change the cell's stored type. Prefixing with a literal
# apostrophe forces Excel to keep/store it as text (quote-prefixed), matching the source.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
}

Set-TextValue $ws.Range("D2") '246.72'
Set-TextValue $ws.Range("D3") '26.43'
Set-TextValue $ws.Range("D4") '5.072'
Set-TextValue $ws.Range("D5") '0.05604'
Set-TextValue $ws.Range("D6") '6.501'
Set-TextValue $ws.Range("D7") '3.049'
Set-TextValue $ws.Range("D9") '0.8404'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range("D10") '0.009819'
$ws.Range("E10").Value = '9OneONEBestin24h'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D11") '0.1345'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D12") '0.02821'
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D13") '0.09386'
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D14") '0.001512'
$ws.Range("E14").Value = '13BitForexTokenBF'
Set-TextValue $ws.Range("D15") '0.006149'
Set-TextValue $ws.Range("D16") '3.555'
Set-TextValue $ws.Range("D18") '0.3181'
Set-TextValue $ws.Range("D19") '0.06962'
Set-TextValue $ws.Range("D20") '0.03115'
Set-TextValue $ws.Range("D21") '0.1301'
Set-TextValue $ws.Range("D22") '3.749'
Set-TextValue $ws.Range("D23") '0.04672'
Set-TextValue $ws.Range("D25") '0.001251'
Set-TextValue $ws.Range("D26") '0.004614'
Set-TextValue $ws.Range("D27") '0.00009594'
Set-TextValue $ws.Range("D28") '0.0001939'
Set-TextValue $ws.Range("D40") '0.03664'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range("D41") '0.006146'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws.Range("D42") '0.1054'
$ws.Range("E42").Value = '41BKEXTokenBKK'
Set-TextValue $ws.Range("D43") '0.002594'
$ws.Range("E43").Value = '42CEJICEJIWorstin24h'
Set-TextValue $ws.Range("D44") '0.008450'
Set-TextValue $ws.Range("D45") '0.00005292'
Set-TextValue $ws.Range("D46") '0.00000000750'
Set-TextValue $ws.Range("D48") '0.002061'
$ws.Range("E48").Value = '47BOLOBOLO'
Set-TextValue $ws.Range("D49") '0.00002099'
Set-TextValue $ws.Range("D50") '0.0001999'
